$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 176, shifting existing rows 176..256 down to 177..257
$ws.Rows.Item(176).Insert()

# Populate the newly inserted row 176 with its full record.
# Columns A,B,C,E,F,G,H,I,N,Q,R are constant across every data row in this sheet;
# D,J,K,L,M,P carry the new record's specific values.
$ws.Cells.Item(176, 1).Value = 8
$ws.Cells.Item(176, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(176, 3).Value = "Coquimbo"
$ws.Cells.Item(176, 4).Value = 45202
$ws.Cells.Item(176, 5).Value = 4
$ws.Cells.Item(176, 6).Value = 100112044
$ws.Cells.Item(176, 7).Value = "Perejil"
$ws.Cells.Item(176, 8).Value = "Sin especificar"
$ws.Cells.Item(176, 9).Value = "Primera"
$ws.Cells.Item(176, 10).Value = 2000
$ws.Cells.Item(176, 11).Value = 1500
$ws.Cells.Item(176, 12).Value = 2000
$ws.Cells.Item(176, 13).Value = 1750
$ws.Cells.Item(176, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(176, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(176, 16).Value = 1167
$ws.Cells.Item(176, 17).Value = 1.5
$ws.Cells.Item(176, 18).Value = "Hortaliza"
